$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate (by text, not a hard-coded index) the empty paragraph that
# immediately follows the "3. He usado ChatGPT para calcular..." paragraph.
# ---------------------------------------------------------------------------
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "3. He usado ChatGPT para calcular*") {
        $anchorIndex = $i
        break
    }
}

$emptyIndex = $anchorIndex + 1

# Create a brand-new, isolated paragraph right after that empty paragraph;
# its Range will be the target for the big InsertXML block below. Re-fetch
# paragraphs by index afterwards since old paragraph handles go stale once
# the collection is mutated.
$d.Paragraphs($emptyIndex).Range.InsertParagraphAfter()
$targetIndex = $emptyIndex + 1
$targetRange = $d.Paragraphs($targetIndex).Range

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">4. He usado </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>chatGPT</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> para resolver la operación (-1) **(</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>a+bE</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>). Me dio la idea de usar números complejos.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">5. </w:t></w:r>
<w:r><w:t xml:space="preserve">He usado </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>chatGPT</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> para hacer más robusta la función de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>calculos</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> de polinomios. Me dio la idea de implementar los </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>checks</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> del tipo de dato que metía</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">if not </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>isinstance</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>coeff</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, list):</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">            raise </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>TypeError</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>"Coefficients must be provided as a list.")</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">        if not </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>coeff</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">            raise </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ValueError</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>"Coefficient list cannot be empty.")</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">        if des not in [0, 1]:</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">            raise </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ValueError</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>"Parameter 'des' must be 0 (ascending) or 1 (descending).")</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
<w:p>
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
</w:body>
</w:document>
'@

$targetRange.InsertXML($xmlFrag)

# ---------------------------------------------------------------------------
# Merge the two trailing (now shifted) empty "Prrafodelista" paragraphs into
# a single paragraph that carries both the pStyle and the en-US language
# run property.
# ---------------------------------------------------------------------------
$listIndices = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "List Paragraph") {
        $listIndices += $i
    }
}

$lastListIndex = $listIndices[$listIndices.Length - 1]
$firstListIndex = $listIndices[$listIndices.Length - 2]

$d.Paragraphs($firstListIndex).Range.Delete()
$d.Paragraphs($firstListIndex).Range.LanguageID = 1033
